$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs at the
# top of the historical list (row 5, right after the 3 "current" summary
# rows), pushing all the existing rows from 5..103 down by one (to 6..104).
# Excel's row Insert() shifts everything down for us, carrying formatting
# along (e.g. the date style on column D), so we just need to insert a row
# and then populate it with the new record's values.

$ws.Rows(5).Insert()

$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 44756
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112024
$ws.Cells.Item(5, 7).Value = "Choclo"
$ws.Cells.Item(5, 8).Value = "Dulce o Americano"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 30000
$ws.Cells.Item(5, 12).Value = 32000
$ws.Cells.Item(5, 13).Value = 31000
$ws.Cells.Item(5, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 443
$ws.Cells.Item(5, 17).Value = 70
$ws.Cells.Item(5, 18).Value = "Hortaliza"
